$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "总计" (so it inherits the
#    same header/index-column style used throughout this workbook), insert
#    it right before "总计", then overwrite its contents with the new
#    7-column per-fund layout (matching the 2021-Qx sheets' structure).
#    NOTE: after Copy(), re-fetch sheets by name -- the original variable
#    can end up bound to the newly inserted copy instead of the source.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Copy($wb.Worksheets.Item("总计"))

$ws = $wb.Worksheets.Item("总计 (2)")
$ws.Name = "2022-Q1"

# Extend the existing header style (currently only on B1:D1) across to H1.
$ws.Range("B1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Header row.
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Only two data rows are needed now (was three on the old "总计" sheet).
$ws.Rows("4:4").Delete()

# Row 2.
$ws.Range("A2").Value = 0
$ws.Range("H2").Value = 3
$ws.Range("B2").Formula = "=""519981"""
$ws.Range("C2").Formula = "=""长信美国标准普尔100等权重指数增强(QDII)"""
$ws.Range("D2").Formula = "=""0.47"""
$ws.Range("E2").Formula = "=""84.16"""
$ws.Range("F2").Formula = "=""0.93"""
$ws.Range("G2").Formula = "=""0.0044"""

# Row 3.
$ws.Range("A3").Value = 1
$ws.Range("H3").Value = 3
$ws.Range("B3").Formula = "=""011706"""
$ws.Range("C3").Formula = "=""长信美国标准普尔100等权重指数增强(QDII) - 美元"""
$ws.Range("D3").Formula = "=""0.47"""
$ws.Range("E3").Formula = "=""84.16"""
$ws.Range("F3").Formula = "=""0.93"""
$ws.Range("G3").Formula = "=""0.0044"""

# Convert the helper formulas above into plain literal text values (so the
# saved cells hold inline text, not formulas) without disturbing styles.
$ws.Range("B2:G3").Copy()
$ws.Range("B2:G3").PasteSpecial(-4163)   # xlPasteValues

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row. Rows("...").Insert() drags the
#    header's bold/border formatting down onto the blank row it creates, so
#    shift the data manually (copy/paste full rows bottom-up) instead.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A4:D4").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4104)   # xlPasteAll
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4:D4").PasteSpecial(-4104)
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4104)

# New row 2.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01

# Re-index column A (0,1,2,3) for the rows that got pushed down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Row 5 was blank before the shift, so it picked up no style on paste;
# restore the A-column index style from row 4 before writing its value.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)      # xlPasteFormats
$totalSheet.Range("A5").Value = 3
